$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell D1 - reuse the header formatting (bold, centered, bordered)
# already applied to A1:C1 by copying its format onto D1.
$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Total Participants"

# Total participants count, placed alongside the first data row
$ws.Range("D2").Value = 3

# New participant rows
$ws.Range("A4").Value = "Bitton"
$ws.Range("B4").Value = "Dan"
$ws.Range("C4").Value = "Oui"

$ws.Range("A5").Value = "Cohen"
$ws.Range("B5").Value = "Yair"
$ws.Range("C5").Value = "Oui"
